# Update cryptocurrency price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "63.842.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  +0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "2.736.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  -0.47%  "
$ws.Range("E4").Value2 = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "565.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "161.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +3.00%  "
$ws.Range("E7").Value2 = "  -0.03%  "
$ws.Range("E8").Value2 = "  -0.63%  "
$ws.Range("E9").Value2 = "  +0.55%  "
$ws.Range("E10").Value2 = "  +4.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "5.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  +1.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.380"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "3.221.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  -0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "26.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +2.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "63.665.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E16").Value2 = "  +0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "2.738.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  -0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "12.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  +3.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "4.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -0.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "355.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "6.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -1.19%  "
$ws.Range("E22").Value2 = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.521"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -2.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "64.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "0.169"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "8.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  +0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "0.0₃0910"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +1.41%  "
$ws.Range("E29").Value2 = "  +3.91%  "
$ws.Range("E30").Value2 = "  +12.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "7.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +2.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "166.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  -1.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "4.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  +1.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "20.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  +0.22%  "
$ws.Range("E35").Value2 = "  +3.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "1.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.978"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  +0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "344.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +5.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "6.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  +2.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "4.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  -0.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "38.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  -0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "21.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  +2.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "21.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  -0.59%  "
$ws.Range("E45").Value2 = "  +0.74%  "
$ws.Range("E46").Value2 = "  +1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.0251"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.0999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "132.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  -1.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.997"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -0.22%  "
$ws.Range("E51").Value2 = "  +0.42%  "
